$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "icon" column (column B) entirely; everything to its right
# shifts left by one column (C->B, D->C, ... H->G).
$ws.Columns(2).Delete()

# Apply left/top alignment to the whole data range (A1:G8) as a single
# style so only one new style entry is created.
$st = $wb.Styles.Add("LeftTop")
$st.HorizontalAlignment = -4131  # xlLeft
$st.VerticalAlignment = -4160    # xlTop
$ws.Range("A1:G8").Style = "LeftTop"

# Update the active selection to match the saved view state.
$ws.Range("C9").Select()

Write-Output "done"
